$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")

$r = $ws1.Range("C1")
$r.Interior.Color = 255
Write-Host "done"
